$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Rows(16).Delete() | Out-Null
$ws.Rows(16).Select() | Out-Null

$lastSheet = $wb.Worksheets.Item("optimization_diagnostics")
$lastSheet.Activate() | Out-Null
